$d = $word.ActiveDocument

# Locate the "Công việc đã làm" progress table (STT / Công việc đã làm),
# rather than hard-coding a table index.
$target = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    $headerText = $candidate.Cell(1, 1).Range.Text
    if ($headerText -match "STT") {
        $target = $candidate
        break
    }
}
if ($target -eq $null) {
    $target = $d.Tables.Item($d.Tables.Count)
}

# Append two new rows at the bottom of the table, mirroring the formatting
# of the existing rows (Rows.Add clones the last row's look automatically).
$newRows = @(
    @{ Stt = "6"; Content = "Kết luận" },
    @{ Stt = "7"; Content = "Hướng phát triển" }
)

foreach ($item in $newRows) {
    $target.Rows.Add() | Out-Null
    $lastRowIndex = $target.Rows.Count
    $target.Cell($lastRowIndex, 1).Range.Text = $item.Stt
    $target.Cell($lastRowIndex, 2).Range.Text = $item.Content
}
